$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted by Excel into a
# number (single "." decimal separator, no thousands grouping). Force them to
# stay plain text (matching the source workbooks inlineStr cells) by briefly
# applying a text number format, then restoring the original "Normal" style so
# the saved XML keeps no explicit s= override (same as before the edit).

$textCells = @("D5","D6","D13","D14","D20","D21","D22","D25","D26","D27","D30","D32","D38","D40","D41","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.839.96"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.463.55"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "575.00"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "145.93"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "2.462.57"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "28.99"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "2.910.83"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "62.747.00"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "2.464.41"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "327.10"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "2.23"
$ws.Range("E22").Value = "  +10.18%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "10.24"
$ws.Range("E25").Value = "  +19.61%  "
$ws.Range("D26").Value = "65.74"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "654.25"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").Value = "0.0₃0982"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "2.584.74"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -13.06%  "
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").Value = "7.99"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "0.370"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "5.38"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "151.03"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "0.0₆0314"
$ws.Range("E44").Value = "  -84.23%  "
$ws.Range("D46").Value = "154.33"
$ws.Range("E46").Value = "  +6.94%  "
$ws.Range("D47").Value = "15.23"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").Value = "3.59"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "20.30"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "0.0512"
$ws.Range("E51").Value = "  -0.04%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
